# TestAddNewTitle.xlsx: the "NewTitle" data value in row 2 (cell C2) was
# changed from "Test_COMPANY_ADMIN_Title" to "RANDOM".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "RANDOM"

# Reflect that C2 is now the selected/active cell, matching the saved
# worksheet's <selection activeCell="C2" sqref="C2"/>.
$ws.Range("C2").Select() | Out-Null
